$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 33333.332
$ws.Range("J3").Value = 33333.332
$ws.Range("L3").Value = 33333.332
$ws.Range("N3").Value = -33561.332
# Row 18
$ws.Range("H18").Value = 1222.9166
$ws.Range("I18").Value = 1152.2727
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 1152.2727
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -868.2727
$ws.Range("N18").Value = -2568
# Row 40
$ws.Range("H40").Value = 1734.1666
$ws.Range("I40").Value = 1552.5
$ws.Range("J40").Value = 1825
$ws.Range("K40").Value = 1552.5
$ws.Range("L40").Value = 1825
$ws.Range("M40").Value = -1377.5
$ws.Range("N40").Value = -2175
# Row 41
$ws.Range("H41").Value = 4982.923
$ws.Range("I41").Value = 162.25
$ws.Range("J41").Value = 7125.4443
$ws.Range("K41").Value = 162.25
$ws.Range("L41").Value = 7125.4443
$ws.Range("M41").Value = 277.75
$ws.Range("N41").Value = -8005.4443
# Row 64
$ws.Range("H64").Value = 2938.182
$ws.Range("I64").Value = 2881.4285
$ws.Range("J64").Value = 3037.5
$ws.Range("K64").Value = 2881.4285
$ws.Range("L64").Value = 3037.5
$ws.Range("M64").Value = -2633.4285
$ws.Range("N64").Value = -3533.5
# Row 67
$ws.Range("H67").Value = 2938.182
$ws.Range("I67").Value = 2881.4285
$ws.Range("J67").Value = 3037.5
$ws.Range("K67").Value = 2881.4285
$ws.Range("L67").Value = 3037.5
$ws.Range("M67").Value = -2023.4285
$ws.Range("N67").Value = -4753.5
# Row 70
$ws.Range("H70").Value = 45052.74
$ws.Range("I70").Value = 111990.336
$ws.Range("J70").Value = 2021.4286
$ws.Range("K70").Value = 335971.008
$ws.Range("L70").Value = 6064.2858
$ws.Range("M70").Value = -335701.008
$ws.Range("N70").Value = -6604.2858
# Row 73
$ws.Range("H73").Value = 45052.74
$ws.Range("I73").Value = 111990.336
$ws.Range("J73").Value = 2021.4286
$ws.Range("K73").Value = 335971.008
$ws.Range("L73").Value = 6064.2858
$ws.Range("M73").Value = -335035.008
$ws.Range("N73").Value = -7936.2858
# Row 75
$ws.Range("H75").Value = 21385.666
$ws.Range("J75").Value = 22078.5
$ws.Range("L75").Value = 22078.5
$ws.Range("N75").Value = -23950.5
# Row 78
$ws.Range("H78").Value = 21385.666
$ws.Range("J78").Value = 22078.5
$ws.Range("L78").Value = 66235.5
$ws.Range("N78").Value = -75595.5
# Row 95
$ws.Range("H95").Value = 47541.332
$ws.Range("J95").Value = 47541.332
$ws.Range("L95").Value = 47541.332
$ws.Range("N95").Value = -53033.332
# Row 102
$ws.Range("H102").Value = 33333.332
$ws.Range("J102").Value = 33333.332
$ws.Range("L102").Value = 33333.332
$ws.Range("N102").Value = -39823.332
# Row 116
$ws.Range("H116").Value = 1678.5714
$ws.Range("I116").Value = 1536.75
$ws.Range("J116").Value = 1867.6666
$ws.Range("K116").Value = 1536.75
$ws.Range("L116").Value = 1867.6666
$ws.Range("M116").Value = 1905.25
$ws.Range("N116").Value = -8751.6666
# Row 129
$ws.Range("H129").Value = 35351.484
$ws.Range("J129").Value = 67871.266
$ws.Range("L129").Value = 203613.798
$ws.Range("N129").Value = -213613.798
# Row 135
$ws.Range("H135").Value = 1240
$ws.Range("I135").Value = 1145
$ws.Range("J135").Value = 1810
$ws.Range("K135").Value = 10305
$ws.Range("L135").Value = 16290
$ws.Range("M135").Value = -7770
$ws.Range("N135").Value = -21360

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1066.174
$ws.Range("I45").Value = 1038.0625
$ws.Range("J45").Value = 1130.4286
$ws.Range("K45").Value = 1038.0625
$ws.Range("L45").Value = 1130.4286
$ws.Range("M45").Value = -661.0625
$ws.Range("N45").Value = -1884.4286
# Row 74
$ws.Range("H74").Value = 1102.6285
$ws.Range("I74").Value = 1047.5625
$ws.Range("J74").Value = 1690
$ws.Range("K74").Value = 1047.5625
$ws.Range("L74").Value = 1690
$ws.Range("M74").Value = -173.5625
$ws.Range("N74").Value = -3438
# Row 77
$ws.Range("H77").Value = 1102.6285
$ws.Range("I77").Value = 1047.5625
$ws.Range("J77").Value = 1690
$ws.Range("K77").Value = 5237.8125
$ws.Range("L77").Value = 8450
$ws.Range("M77").Value = -869.8125
$ws.Range("N77").Value = -17186
# Row 110
$ws.Range("H110").Value = 5674.36
$ws.Range("I110").Value = 7120.353
$ws.Range("K110").Value = 7120.353
$ws.Range("M110").Value = -5075.353

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 421.53845
$ws.Range("I22").Value = 421.53845
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 421.53845
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -248.53845
$ws.Range("N22").ClearContents()
# Row 76
$ws.Range("H76").Value = 25078.5
$ws.Range("J76").Value = 30104.666
$ws.Range("L76").Value = 30104.666
$ws.Range("N76").Value = -30734.666
# Row 79
$ws.Range("H79").Value = 25078.5
$ws.Range("J79").Value = 30104.666
$ws.Range("L79").Value = 30104.666
$ws.Range("N79").Value = -32288.666
# Row 92
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
# Row 128
$ws.Range("H128").Value = 1610
$ws.Range("I128").Value = 1610
$ws.Range("K128").Value = 4830
$ws.Range("M128").Value = -2340

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 43
$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20368
# Row 58
$ws.Range("H58").Value = 1184.6875
$ws.Range("I58").Value = 1012.5
$ws.Range("J58").Value = 1471.6666
$ws.Range("K58").Value = 1012.5
$ws.Range("L58").Value = 1471.6666
$ws.Range("M58").Value = -809.5
$ws.Range("N58").Value = -1877.6666
# Row 101
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490
# Row 136
$ws.Range("H136").Value = 1184.6875
$ws.Range("I136").Value = 1012.5
$ws.Range("J136").Value = 1471.6666
$ws.Range("K136").Value = 3037.5
$ws.Range("L136").Value = 4414.9998
$ws.Range("M136").Value = -487.5
$ws.Range("N136").Value = -9514.9998

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 4684.6665
$ws.Range("I131").Value = 7702.6
$ws.Range("J131").Value = 912.25
$ws.Range("K131").Value = 23107.8
$ws.Range("L131").Value = 2736.75
$ws.Range("M131").Value = -18067.8
$ws.Range("N131").Value = -12816.75

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1377.6666
$ws.Range("I102").Value = 1242.7142
$ws.Range("K102").Value = 1242.7142
$ws.Range("M102").Value = 379.2858000000001
# Row 126
$ws.Range("H126").Value = 2091.2
$ws.Range("I126").Value = 1940.6666
$ws.Range("J126").Value = 2542.8
$ws.Range("K126").Value = 5821.9998
$ws.Range("L126").Value = 7628.400000000001
$ws.Range("M126").Value = -3351.9998
$ws.Range("N126").Value = -12568.4

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 781.8182
$ws.Range("I100").Value = 962.5
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = 962.5
$ws.Range("L100").Value = 300
$ws.Range("M100").Value = -421.5
$ws.Range("N100").Value = -1382
# Row 132
$ws.Range("H132").Value = 2256.8635
$ws.Range("I132").Value = 2458
$ws.Range("J132").Value = 2055.7273
$ws.Range("K132").Value = 7374
$ws.Range("L132").Value = 6167.1819
$ws.Range("M132").Value = -4844
$ws.Range("N132").Value = -11227.1819

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 22388.5
$ws.Range("J80").Value = 22388.5
$ws.Range("L80").Value = 22388.5
$ws.Range("N80").Value = -24384.5
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 83
$ws.Range("H83").Value = 22388.5
$ws.Range("J83").Value = 22388.5
$ws.Range("L83").Value = 67165.5
$ws.Range("N83").Value = -77149.5
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 97
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982
# Row 101
$ws.Range("H101").Value = 20038.4
$ws.Range("J101").Value = 20038.4
$ws.Range("L101").Value = 20038.4
$ws.Range("N101").Value = -26528.4
# Row 122
$ws.Range("H122").Value = 52275.4
$ws.Range("I122").Value = 79122.16
$ws.Range("J122").Value = 2417.1428
$ws.Range("K122").Value = 237366.48
$ws.Range("L122").Value = 7251.428400000001
$ws.Range("M122").Value = -234916.48
$ws.Range("N122").Value = -12151.4284
